# Add new daily rows (419-449) to each of the four province sheets.
# Rows 419-421 carry new case/death/recovery/hospitalisation counts (column C)
# together with the rolling 7-day average in column D (continuing the
# existing shared formula pattern). Rows 422-449 only carry the date in
# column A (data not yet available for those days).

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{ Name = "Nuovi casi";        C419 = 37; C420 = 27; C421 = 2  },
    @{ Name = "Deceduti";          C419 = 1;  C420 = 2;  C421 = 1  },
    @{ Name = "Dimessi   Guariti"; C419 = 28; C420 = 3;  C421 = 63 },
    @{ Name = "Ricoveri";          C419 = 35; C420 = 37; C421 = 35 }
)

$firstNewDate = 44317   # 2021-05-01 serial date for row 419
$lastDateRow  = 449     # last row to receive a date-only entry

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Name)

    # --- Rows 419:421 -- new observations -----------------------------
    $ws.Cells.Item(419, 1).Value = $firstNewDate
    $ws.Cells.Item(420, 1).Value = $firstNewDate + 1
    $ws.Cells.Item(421, 1).Value = $firstNewDate + 2

    $ws.Cells.Item(419, 3).Value = $sd.C419
    $ws.Cells.Item(420, 3).Value = $sd.C420
    $ws.Cells.Item(421, 3).Value = $sd.C421

    # Match formatting of the preceding data row, then extend the
    # rolling 7-day AVERAGE formula (column D) down through row 421.
    $ws.Range("D418").Copy()
    $ws.Range("D419:D421").PasteSpecial(-4122)
    $ws.Range("D419:D421").Formula = "=AVERAGE(C413:C419)"

    # --- Rows 422:449 -- dates only (no counts yet) --------------------
    for ($r = 422; $r -le $lastDateRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $firstNewDate + ($r - 419)
    }

    # Selection mirrors what was left behind after pasting the new block.
    $ws.Range("C419:C421").Select()
}

# The last sheet edited ("Ricoveri") ends up the active one, with the
# cursor resting on the newest rolling-average cell.
$wsLast = $wb.Worksheets.Item("Ricoveri")
$wsLast.Range("D421").Select()
